$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header-row labels (A1:U1) so the "_old"/"_new" suffixes used
#    for the two compared format versions become the concrete version tags
#    "_FV2310" (previous) / "_FV2404" (current). The "diff" header (K1) is
#    left untouched.
# ---------------------------------------------------------------------------
$headerNames = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headerNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerNames[$i]
}

# ---------------------------------------------------------------------------
# 2. Freeze the header row: top-left pane stays put, the row-2-down pane is
#    frozen so the header (row 1) remains visible while scrolling.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Turn the used range A1:U93 into a native Excel table ("Table1") with an
#    autofilter on the header row, matching the already-present column
#    headers, and no explicit table style (keep workbook default look).
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U93")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Leave the selection on the header cell, same as the workbook's normal
# "just opened" resting state.
$ws.Range("A1").Select()
